$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" '25.953.69'
Set-TextValue "E2" '  -0.33%  '
Set-TextValue "D3" '1.621.60'
Set-TextValue "E3" '  -0.95%  '
Set-TextValue "E4" '  -0.17%  '
Set-TextValue "D5" '212.71'
Set-TextValue "E5" '  -1.00%  '
Set-TextValue "E6" '  -0.15%  '
Set-TextValue "D7" '0.488'
Set-TextValue "E7" '  -3.40%  '
Set-TextValue "D8" '0.0622'
Set-TextValue "E8" '  -0.97%  '
Set-TextValue "E9" '  -1.30%  '
Set-TextValue "D10" '18.33'
Set-TextValue "E10" '  -2.09%  '
Set-TextValue "E11" '  -0.21%  '
Set-TextValue "D12" '1.846.49'
Set-TextValue "E12" '  -1.00%  '
Set-TextValue "D13" '1.614.08'
Set-TextValue "E13" '  -3.49%  '
Set-TextValue "E14" '  -1.55%  '
Set-TextValue "E15" '  -1.85%  '
Set-TextValue "D16" '25.960.12'
Set-TextValue "E16" '  -0.35%  '
Set-TextValue "D17" '61.71'
Set-TextValue "E17" '  -0.92%  '
Set-TextValue "E18" '  -1.36%  '
Set-TextValue "D19" '1.01'
Set-TextValue "E19" '  -0.21%  '
Set-TextValue "D20" '191.94'
Set-TextValue "E20" '  +0.31%  '
Set-TextValue "E21" '  -0.63%  '
Set-TextValue "E22" '  -0.93%  '
Set-TextValue "E23" '  -2.23%  '
Set-TextValue "E24" '  -0.12%  '
Set-TextValue "D25" '144.35'
Set-TextValue "E26" '  -0.19%  '
Set-TextValue "D27" '1.71'
Set-TextValue "E27" '  -3.75%  '
Set-TextValue "E28" '  -1.91%  '
Set-TextValue "D29" '15.20'
Set-TextValue "E29" '  -0.58%  '
Set-TextValue "E30" '  -1.08%  '
Set-TextValue "E31" '  -1.58%  '
Set-TextValue "E32" '  -1.76%  '
Set-TextValue "E33" '  -2.83%  '
Set-TextValue "E34" '  -0.53%  '
Set-TextValue "E35" '  -1.39%  '
Set-TextValue "D36" '1.127.97'
Set-TextValue "E36" '  -0.40%  '
Set-TextValue "E37" '  -4.02%  '
Set-TextValue "E38" '  -1.92%  '
Set-TextValue "E39" '  -2.11%  '
Set-TextValue "E40" '  -1.56%  '
Set-TextValue "D41" '97.96'
Set-TextValue "E41" '  -1.01%  '
Set-TextValue "D42" '1.757.77'
Set-TextValue "E42" '  -0.88%  '
Set-TextValue "E43" '  -4.12%  '
Set-TextValue "E44" '  -3.32%  '
Set-TextValue "E45" '  -0.47%  '
Set-TextValue "E46" '  +1.48%  '
Set-TextValue "D47" '54.06'
Set-TextValue "E47" '  -2.77%  '
Set-TextValue "D48" '0.0516'
Set-TextValue "E48" '  -1.64%  '
Set-TextValue "E49" '  -0.83%  '
Set-TextValue "D50" '7.48'
Set-TextValue "E50" '  -1.63%  '
Set-TextValue "E51" '  +0.03%  '
